$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS group)
$ws.Range("D2").Value = 906.49
$ws.Range("E2").Value = -906.49

# Row 4 (TOTAL row)
$ws.Range("D4").Value = 1458.92
$ws.Range("E4").Value = 12264.42
$ws.Range("F4").Value = 0.1063093969835332
